# Insert a new data row before row 121 (this pushes the existing rows
# 121:193 down to 122:194, preserving all of their data/formatting).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("121:121").Insert()

# Populate the newly inserted row 121 with the new record.
$ws.Range("A121").Value = 10
$ws.Range("B121").Value = "Vega Modelo de Temuco"
$ws.Range("C121").Value = "La Araucanía"
$ws.Range("D121").Value = 44438
$ws.Range("E121").Value = 9
$ws.Range("F121").Value = "Fruta"
$ws.Range("G121").Value = 100108
$ws.Range("H121").Value = "Tropicales y subtropicales"
$ws.Range("I121").Value = 100108002
$ws.Range("J121").Value = "Mango"
$ws.Range("K121").Value = "Sin especificar"
$ws.Range("L121").Value = "Primera"
$ws.Range("M121").Value = 380
$ws.Range("N121").Value = 9000
$ws.Range("O121").Value = 9000
$ws.Range("P121").Value = 9000
$ws.Range("Q121").Value = "$/bandeja 4 kilos"
$ws.Range("R121").Value = "Brasil"
$ws.Range("S121").Value = 2250
$ws.Range("T121").Value = 4

# Make sure the reported dimension/used range covers the new row count.
$ws.Range("A1:T194").Select()
